$wb = $excel.ActiveWorkbook

# ----- Sheet 1: "VENTAS POR GRUPO" -----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L5").Value = 940.62
$ws1.Range("M5").Value = 6733.24
$ws1.Range("D11").Value = 648
$ws1.Range("M11").Value = 6208.03
$ws1.Range("M13").Value = -494.21
$ws1.Range("D22").Value = 734.4299999999999
$ws1.Range("H22").Value = 1858.63
$ws1.Range("L23").Value = "2 de 21"

# ----- Sheet 2: "VENTA MENSUAL" -----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F5").Value = 7673.86
$ws2.Range("F11").Value = 6856.03
$ws2.Range("F13").Value = -494.21
$ws2.Range("F22").Value = 4205.02
$ws2.Range("F23").Value = 27396.16
# raw OOXML column width of 13 corresponds to a COM ColumnWidth of (13 - 5/6)
$ws2.Columns.Item(6).ColumnWidth = 13 - 5/6

# ----- Sheet 3: "CUMPLIMIENTO MENSUAL" -----
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 1382.43
$ws3.Range("E3").Value = 4122.18890386263
$ws3.Range("F3").Value = 0.2511400015412401

$ws3.Range("D6").Value = 1858.63
$ws3.Range("E6").Value = 1048.95368146026
$ws3.Range("F6").Value = 0.6392352563577983

$ws3.Range("D11").Value = 2552.58
$ws3.Range("E11").Value = 3291.86916370549
$ws3.Range("F11").Value = 0.4367528792707672

$ws3.Range("D12").Value = 21602.52
$ws3.Range("E12").Value = 16137.22
$ws3.Range("F12").Value = 0.5724077590359659

$ws3.Range("D14").Value = 27396.16
$ws3.Range("E14").Value = 28028.58147880389
$ws3.Range("F14").Value = 0.4942947728583836

# raw OOXML column width of 22 corresponds to a COM ColumnWidth of (22 - 5/6)
$ws3.Columns.Item(5).ColumnWidth = 22 - 5/6

Write-Host "applied changes"
